# Auto-generated script to refresh market-price derived columns (H-N)
# on the Zodiark_Profits leve-crafting sheets, per the scheduled data-refresh run.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2762
$ws.Range("J17").Value = 2762
$ws.Range("L17").Value = 8286
$ws.Range("N17").Value = -8622

$ws.Range("H39").Value = 637.5238000000001
$ws.Range("I39").Value = 77.3
$ws.Range("J39").Value = 1146.8182
$ws.Range("K39").Value = 231.9
$ws.Range("L39").Value = 3440.4546
$ws.Range("M39").Value = 64.10000000000002
$ws.Range("N39").Value = -4032.4546

$ws.Range("H58").Value = 11064.286
$ws.Range("I58").Value = 950
$ws.Range("K58").Value = 2850
$ws.Range("M58").Value = -2700

$ws.Range("H62").Value = 4093.2856
$ws.Range("I62").Value = 4076.5
$ws.Range("K62").Value = 4076.5
$ws.Range("M62").Value = -3452.5

$ws.Range("H65").Value = 4093.2856
$ws.Range("I65").Value = 4076.5
$ws.Range("K65").Value = 20382.5
$ws.Range("M65").Value = -17262.5

$ws.Range("H86").Value = 4435.577
$ws.Range("I86").Value = 4986.357
$ws.Range("J86").Value = 3793
$ws.Range("K86").Value = 4986.357
$ws.Range("L86").Value = 3793
$ws.Range("M86").Value = -3863.357
$ws.Range("N86").Value = -6039

$ws.Range("H89").Value = 4435.577
$ws.Range("I89").Value = 4986.357
$ws.Range("J89").Value = 3793
$ws.Range("K89").Value = 24931.785
$ws.Range("L89").Value = 18965
$ws.Range("M89").Value = -19315.785
$ws.Range("N89").Value = -30197

$ws.Range("H137").Value = 1689.6957
$ws.Range("I137").Value = 1524.579
$ws.Range("J137").Value = 2474
$ws.Range("K137").Value = 4573.737
$ws.Range("L137").Value = 7422
$ws.Range("M137").Value = -2023.737
$ws.Range("N137").Value = -12522

$ws.Range("H138").Value = 2939.5417
$ws.Range("I138").Value = 1400.2693
$ws.Range("J138").Value = 3809.5652
$ws.Range("K138").Value = 4200.8079
$ws.Range("L138").Value = 11428.6956
$ws.Range("M138").Value = 939.1921000000002
$ws.Range("N138").Value = -21708.6956

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4491.731
$ws.Range("I32").Value = 2898.186
$ws.Range("K32").Value = 2898.186
$ws.Range("M32").Value = -2611.186

$ws.Range("H45").Value = 1495.4138
$ws.Range("I45").Value = 1491.9375
$ws.Range("K45").Value = 1491.9375
$ws.Range("M45").Value = -1114.9375

$ws.Range("H61").Value = 2949.8857
$ws.Range("I61").Value = 2419.7144
$ws.Range("K61").Value = 2419.7144
$ws.Range("M61").Value = -2207.7144

$ws.Range("H74").Value = 3266.9
$ws.Range("I74").Value = 3228.8235
$ws.Range("J74").Value = 3482.6667
$ws.Range("K74").Value = 3228.8235
$ws.Range("L74").Value = 3482.6667
$ws.Range("M74").Value = -2354.8235
$ws.Range("N74").Value = -5230.6667

$ws.Range("H77").Value = 3266.9
$ws.Range("I77").Value = 3228.8235
$ws.Range("J77").Value = 3482.6667
$ws.Range("K77").Value = 16144.1175
$ws.Range("L77").Value = 17413.3335
$ws.Range("M77").Value = -11776.1175
$ws.Range("N77").Value = -26149.3335

$ws.Range("H110").Value = 573.4545000000001
$ws.Range("I110").Value = 573.4545000000001
$ws.Range("K110").Value = 573.4545000000001
$ws.Range("M110").Value = 1471.5455

$ws.Range("H122").Value = 8888
$ws.Range("I122").Value = 8888
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 26664
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -24214
$ws.Range("N122").ClearContents()

$ws.Range("H132").Value = 11529.333
$ws.Range("I132").Value = 8118.985
$ws.Range("J132").Value = 40091
$ws.Range("K132").Value = 24356.955
$ws.Range("L132").Value = 120273
$ws.Range("M132").Value = -21826.955
$ws.Range("N132").Value = -125333

$ws.Range("H136").Value = 2949.8857
$ws.Range("I136").Value = 2419.7144
$ws.Range("K136").Value = 7259.1432
$ws.Range("M136").Value = -4709.1432

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3746.4783
$ws.Range("I20").Value = 2826.5881
$ws.Range("J20").Value = 6352.8335
$ws.Range("K20").Value = 2826.5881
$ws.Range("L20").Value = 6352.8335
$ws.Range("M20").Value = -2579.5881
$ws.Range("N20").Value = -6846.8335

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H23").Value = 45005
$ws.Range("I23").Value = 0
$ws.Range("J23").Value = 45005
$ws.Range("K23").Value = 0
$ws.Range("L23").Value = 45005
$ws.Range("N23").Value = -45485
$ws.Range("M23").ClearContents()

$ws.Range("H27").Value = 45005
$ws.Range("I27").Value = 0
$ws.Range("J27").Value = 45005
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = 45005
$ws.Range("N27").Value = -45389
$ws.Range("M27").ClearContents()

$ws.Range("H31").Value = 1942.25
$ws.Range("I31").Value = 2229.8333
$ws.Range("J31").Value = 1819
$ws.Range("K31").Value = 2229.8333
$ws.Range("L31").Value = 1819
$ws.Range("M31").Value = -1934.8333
$ws.Range("N31").Value = -2409

$ws.Range("H34").Value = 1942.25
$ws.Range("I34").Value = 2229.8333
$ws.Range("J34").Value = 1819
$ws.Range("K34").Value = 2229.8333
$ws.Range("L34").Value = 1819
$ws.Range("M34").Value = -2027.8333
$ws.Range("N34").Value = -2223

$ws.Range("H58").Value = 4820.9
$ws.Range("I58").Value = 4890.421
$ws.Range("J58").Value = 3500
$ws.Range("K58").Value = 4890.421
$ws.Range("L58").Value = 3500
$ws.Range("M58").Value = -4687.421
$ws.Range("N58").Value = -3906

$ws.Range("H132").Value = 1617.5555
$ws.Range("I132").Value = 1617.5555
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 4852.666499999999
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -2322.666499999999
$ws.Range("N132").ClearContents()

$ws.Range("H136").Value = 4820.9
$ws.Range("I136").Value = 4890.421
$ws.Range("J136").Value = 3500
$ws.Range("K136").Value = 14671.263
$ws.Range("L136").Value = 10500
$ws.Range("M136").Value = -12121.263
$ws.Range("N136").Value = -15600

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 566.6667
$ws.Range("I5").Value = 566.6667
$ws.Range("K5").Value = 1700.0001
$ws.Range("M5").Value = -1588.0001

$ws.Range("H21").Value = 0
$ws.Range("I21").Value = 0
$ws.Range("K21").Value = 0
$ws.Range("M21").ClearContents()

$ws.Range("H131").Value = 1533.6
$ws.Range("J131").Value = 2064.5
$ws.Range("L131").Value = 6193.5
$ws.Range("N131").Value = -16273.5

$ws.Range("H135").Value = 566.6667
$ws.Range("I135").Value = 566.6667
$ws.Range("K135").Value = 5100.0003
$ws.Range("M135").Value = -2565.0003

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 28159.678
$ws.Range("I70").Value = 50374.332
$ws.Range("J70").Value = 11498.6875
$ws.Range("K70").Value = 50374.332
$ws.Range("L70").Value = 11498.6875
$ws.Range("M70").Value = -50104.332
$ws.Range("N70").Value = -12038.6875

$ws.Range("H73").Value = 28159.678
$ws.Range("I73").Value = 50374.332
$ws.Range("J73").Value = 11498.6875
$ws.Range("K73").Value = 50374.332
$ws.Range("L73").Value = 11498.6875
$ws.Range("M73").Value = -49438.332
$ws.Range("N73").Value = -13370.6875

$ws.Range("H113").Value = 1139.8
$ws.Range("I113").Value = 1084
$ws.Range("J113").Value = 1223.5
$ws.Range("K113").Value = 1084
$ws.Range("L113").Value = 1223.5
$ws.Range("M113").Value = 1086
$ws.Range("N113").Value = -5563.5

$ws.Range("H126").Value = 3928.0715
$ws.Range("I126").Value = 3291.9092
$ws.Range("K126").Value = 9875.7276
$ws.Range("M126").Value = -7405.7276

$ws.Range("H132").Value = 2138.15
$ws.Range("I132").Value = 828.6923
$ws.Range("K132").Value = 2486.0769
$ws.Range("M132").Value = 43.92309999999998

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 8710
$ws.Range("I16").Value = 8710
$ws.Range("K16").Value = 8710
$ws.Range("M16").Value = -8540

$ws.Range("H93").Value = 3395.7368
$ws.Range("I93").Value = 1885.9286
$ws.Range("K93").Value = 1885.9286
$ws.Range("M93").Value = -637.9286

$ws.Range("H132").Value = 5839.125
$ws.Range("I132").Value = 5349.6665
$ws.Range("J132").Value = 7307.5
$ws.Range("K132").Value = 16048.9995
$ws.Range("L132").Value = 21922.5
$ws.Range("M132").Value = -13518.9995
$ws.Range("N132").Value = -26982.5

$ws.Range("H136").Value = 6644.207
$ws.Range("I136").Value = 5642.25
$ws.Range("K136").Value = 16926.75
$ws.Range("M136").Value = -14376.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 2666.1428
$ws.Range("I126").Value = 2399.6
$ws.Range("J126").Value = 3332.5
$ws.Range("K126").Value = 7198.799999999999
$ws.Range("L126").Value = 9997.5
$ws.Range("M126").Value = -4728.799999999999
$ws.Range("N126").Value = -14937.5

$ws.Range("H132").Value = 1630.6333
$ws.Range("I132").Value = 1585.5769
$ws.Range("J132").Value = 1923.5
$ws.Range("K132").Value = 4756.7307
$ws.Range("L132").Value = 5770.5
$ws.Range("M132").Value = -2226.7307
$ws.Range("N132").Value = -10830.5

Write-Host "Applied Zodiark_Profits market data refresh"
